# ----------------------------------------------------------------------------
# Applies the "add ground_warfare files and data" edit:
#  - "ship" sheet: extend data table from columns B:E to B:G (6 data columns
#    instead of 4), update several numeric values, and strip the red/yellow
#    highlight fills that used to flag mismatched cells (now all clean).
#  - "SSM" sheet: update speed (row2) and attack_range (row8) figures.
#  - "SAM" sheet: update speed (row2) and attack_range (row3) figures.
#  - "inception" sheet: update inception_distance (row2) and widen column A.
#  - refresh the active selections left behind by the editing session.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "ship" sheet — rebuild the B:G data block and clear the old fills.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ship")

$shipData = @{}
$shipData[1]  = @(1,2,3,4,5,6)
$shipData[2]  = @(25,25,25,25,25,25)
$shipData[3]  = @(90,90,90,90,90,90)
$shipData[4]  = @(3,3,3,3,3,3)
$shipData[5]  = @(12,10,10,10,10,10)
$shipData[6]  = @(22,24,24,24,24,24)
$shipData[7]  = @(36,36,36,36,36,36)
$shipData[8]  = @(20,20,20,20,20,20)
$shipData[9]  = @(3,3,3,3,3,3)
$shipData[10] = @(15,12,12,12,12,12)
$shipData[11] = @(1,1,1,1,1,1)
$shipData[12] = @(15,20,20,20,20,20)
$shipData[13] = @(1,1,1,1,1,1)
$shipData[14] = @(1,3,3,4,5,3)
$shipData[15] = @(15,8,8,8,8,8)
$shipData[16] = @(2,2,2,2,2,2)
$shipData[17] = @(20,20,20,20,20,20)
$shipData[18] = @(-0.1,-0.1,-0.1,-0.1,-0.1,-0.1)
$shipData[19] = @(10,10,10,10,10,10)
$shipData[20] = @(0.9,0.9,0.9,0.9,0.9,0.9)
$shipData[21] = @(100,100,100,100,100,100)
$shipData[22] = @(25,25,25,25,25,25)
$shipData[23] = @(20,20,20,20,20,20)
$shipData[24] = @(150,150,150,150,150,150)
$shipData[25] = @(6000000,6000000,6000000,6000000,6000000,6000000)
$shipData[26] = @(9300,9300,9300,9300,9300,9300)
$shipData[27] = @(7.5,7.5,7.5,7.5,7.5,7.5)
$shipData[28] = @(4,4,4,4,4,4)
$shipData[29] = @(4,4,4,4,4,4)
$shipData[30] = @(4200,4200,4200,4200,4200,4200)
$shipData[31] = @(1100,1100,1100,1100,1100,1100)
$shipData[32] = @(10,10,10,10,10,10)
$shipData[33] = @(11,11,11,11,11,11)
$shipData[34] = @(2,2,2,2,2,2)
$shipData[35] = @(3,3,3,3,3,3)
$shipData[36] = @(2,2,2,2,2,2)
$shipData[37] = @(3,3,3,3,3,3)

foreach ($r in $shipData.Keys) {
    $vals = $shipData[$r]
    for ($ci = 0; $ci -lt $vals.Length; $ci++) {
        $ws1.Cells.Item($r, 2 + $ci).Value = $vals[$ci]
    }
}

# row 38 holds the text "yellow" flag in every data column (B..G)
$ws1.Cells.Item(38, 2).Value = "blue"
$ws1.Cells.Item(38, 3).Value = "yellow"
$ws1.Cells.Item(38, 4).Value = "yellow"
$ws1.Cells.Item(38, 5).Value = "yellow"
$ws1.Cells.Item(38, 6).Value = "yellow"
$ws1.Cells.Item(38, 7).Value = "yellow"

# the old red/yellow conditional fills (style indexes 1 & 2) are gone —
# clear all cell formatting on the sheet back to the default style and
# drop the per-column styled widths.
$ws1.Cells.ClearFormats()

$ws1.Activate()
$ws1.Range("G14").Select()

# ---------------------------------------------------------------------------
# 2. "SSM" sheet — speed & attack_range updates.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("SSM")
$ws2.Range("B2").Value = 2.5
$ws2.Range("C2").Value = 3
$ws2.Range("D2").Value = 4.5
$ws2.Range("B8").Value = 150
$ws2.Range("C8").Value = 150
$ws2.Range("D8").Value = 150

$ws2.Activate()
$ws2.Range("C3").Select()

# ---------------------------------------------------------------------------
# 3. "SAM" sheet — speed & attack_range updates.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("SAM")
$ws3.Range("B2").Value = 5.5
$ws3.Range("C2").Value = 5.5
$ws3.Range("D2").Value = 5.5
$ws3.Range("E2").Value = 5.5
$ws3.Range("B3").Value = 60
$ws3.Range("C3").Value = 60
$ws3.Range("D3").Value = 30
$ws3.Range("E3").Value = 30

$ws3.Activate()
$ws3.Range("N11").Select()

# ---------------------------------------------------------------------------
# 4. "inception" sheet — inception_distance update + wider id column.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("inception")
$ws5.Range("B2").Value = 90

$ws5.Columns("A").ColumnWidth = 15.25

$ws5.Activate()
$ws5.Range("B2").Select()

# ---------------------------------------------------------------------------
# Leave "ship" as the active/visible sheet, matching tabSelected="1".
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("G14").Select()
